$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190981030464172
$ws.Range("B1").Value = 2.579953193664551
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.199592351913452
$ws.Range("E1").Value = 1.180020093917847
